$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the cell to be stored as text (shared string) rather than a
    # number, then strip the temporary "Text" number-format back off so the
    # cell's style stays at its original (default) index.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- Sheet "RTECreation": update RWTrackingNo values in column C ---
$wsCreation = $wb.Worksheets.Item("RTECreation")
Set-TextValue $wsCreation.Range("C2") "125980979"
Set-TextValue $wsCreation.Range("C3") "125980991"

# --- Sheet "SearchRTE": update RWTrackingNo / JobID / PickUPID / BOLNo values ---
$wsSearch = $wb.Worksheets.Item("SearchRTE")
Set-TextValue $wsSearch.Range("A2") "125980979"
Set-TextValue $wsSearch.Range("B2") "32393886"
Set-TextValue $wsSearch.Range("C2") "3399124"
Set-TextValue $wsSearch.Range("D2") "125980980"

Set-TextValue $wsSearch.Range("A3") "125980991"
Set-TextValue $wsSearch.Range("B3") "32393887"
Set-TextValue $wsSearch.Range("C3") "3399125"
Set-TextValue $wsSearch.Range("D3") "125981002"
